$wb = $excel.ActiveWorkbook

# "Generate Report for Handback" - refresh the handoff/handback timestamps
# for the 55674320-.../f22db898-... file pair on both the zh-cn and de-de
# report sheets (rows 26 and 27 share the same Correspond Handoff/Handback
# datetimes).

$ws_zh = $wb.Worksheets.Item("zh-cn")
$ws_de = $wb.Worksheets.Item("de-de")

# zh-cn: Correspond Handoff Datetime (col D) / Correspond Handback DateTime (col G)
$ws_zh.Range("D26").Value = "2016-03-03 08:15:04"
$ws_zh.Range("D27").Value = "2016-03-03 08:15:04"
$ws_zh.Range("G26").Value = "2016-03-03 08:15:54"
$ws_zh.Range("G27").Value = "2016-03-03 08:15:54"

# de-de: Correspond Handoff Datetime (col D) / Correspond Handback DateTime (col G)
$ws_de.Range("D26").Value = "2016-03-03 08:15:15"
$ws_de.Range("D27").Value = "2016-03-03 08:15:15"
$ws_de.Range("G26").Value = "2016-03-03 08:16:16"
$ws_de.Range("G27").Value = "2016-03-03 08:16:16"
